$wb = $excel.ActiveWorkbook

# Source sheet to clone the layout/styles from
$denmark = $wb.Worksheets.Item("Denmark")

# ---------------------------------------------------------------------
# Russia (copy of Denmark, then drop the MZXSDR240 row -> 19 rows)
# ---------------------------------------------------------------------
$denmark.Copy($null, $denmark)
$russia = $wb.Worksheets.Item($denmark.Index + 1)
$russia.Name = "Russia"
$russia.Range("B4").Value = "NGC-2929/T2923  "
$russia.Range("B2").Value = "Russia Market"
$russia.Rows("16").Delete()
$russia.Range("A1:D19").Select()

# ---------------------------------------------------------------------
# Finland (copy of Denmark, keeps all 20 rows unchanged)
# ---------------------------------------------------------------------
$denmark.Copy($null, $russia)
$finland = $wb.Worksheets.Item($russia.Index + 1)
$finland.Name = "Finland"
$finland.Range("B4").Value = "NGC-3130/T2886  "
$finland.Range("B2").Value = "Finland Market"
$finland.Range("A1:D19").Select()

# ---------------------------------------------------------------------
# Hungary (copy of Russia layout - 19 rows, becomes the active tab)
# ---------------------------------------------------------------------
$russia.Copy($null, $finland)
$hungary = $wb.Worksheets.Item($finland.Index + 1)
$hungary.Name = "Hungary"
$hungary.Range("B4").Value = "NGC-3104/T3004  "
$hungary.Range("B2").Value = "Hungary Market"
$hungary.Range("A1:D19").Select()

$hungary.Activate()
